$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.744.52'
$ws.Range("E2").Value = '  +0.28%  '

$ws.Range("D3").Value = '1.649.51'
$ws.Range("E3").Value = '  +0.72%  '

$ws.Range("E4").Value = '  +0.19%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.83'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.22%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.504'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.81%  '

$ws.Range("E7").Value = '  +0.21%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.253'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.64%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0628'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.68%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.36'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.80%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0844'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.32%  '

$ws.Range("E12").Value = '  +0.85%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.681.49'
$ws.Range("E13").Value = '  +2.79%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.22'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.97%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.534'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.21%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.69'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.16%  '

$ws.Range("D17").Value = '26.828.37'
$ws.Range("E17").Value = '  +0.60%  '

$ws.Range("D18").Value = '0.0₃0755'
$ws.Range("E18").Value = '  +1.32%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '220.88'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.99%  '

$ws.Range("E20").Value = '  +0.13%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.40'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.05%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.35'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.12%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.56'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.67%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.13'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +11.09%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.88'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.23%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.01'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.56%  '

$ws.Range("E27").Value = '  +0.18%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.10'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.61%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.91'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.16%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0519'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.15%  '

$ws.Range("E31").Value = '  +0.65%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.42'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.27%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.05'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.49%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.56'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.20%  '

$ws.Range("D35").Value = '1.290.43'
$ws.Range("E35").Value = '  +7.90%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0182'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.80%  '

$ws.Range("E37").Value = '  +1.20%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.831'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.57%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.525'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.94%  '

$ws.Range("E40").Value = '  +0.19%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.813'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.32%  '

$ws.Range("E42").Value = '  -2.75%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.45'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.73%  '

$ws.Range("D44").Value = '1.792.29'
$ws.Range("E44").Value = '  +1.04%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '93.62'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.54%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '59.59'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +8.68%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.61'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.68%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0517'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.99%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.79'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.58%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0976'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.85%  '

$ws.Range("E51").Value = '  -0.41%  '
